$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INIT BOM")

# --- Row 8: MCU ---
$ws.Range("B8").Value = "MCU"
$ws.Range("C8").Value = "STM32H742VGT6"
$ws.Range("D8").Value = "ST"
$ws.Range("E8").NumberFormat = "0.00"
$ws.Range("E8").Value = 15.8147
$ws.Range("F8").Value = "https://www.st.com/resource/en/datasheet/stm32h743zg.pdf"
$ws.Range("H8").Value = "https://www.digikey.nl/nl/products/detail/stmicroelectronics/STM32H742VGT6/12337748"

# --- Row 9: POWER INPUT JACK ---
$ws.Range("B9").Value = "POWER INPUT JACK"
$ws.Range("D9").Value = "WURTH"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "694106301002"
$ws.Range("E9").NumberFormat = "0.00"
$ws.Range("E9").Value = 1.1374
$ws.Range("F9").Value = "https://www.we-online.com/katalog/datasheet/6941xx301002.pdf"
$ws.Range("G9").Value = "https://www.digikey.be/en/products/detail/w%C3%BCrth-elektronik/694106301002/5047522"

# --- Hyperlinks for the new datasheet cells (re-apply the Hyperlink cell
# style afterwards so it reuses the workbook's existing Hyperlink style
# rather than minting a brand-new one) ---
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.st.com/resource/en/datasheet/stm32h743zg.pdf") | Out-Null
$ws.Range("F8").Style = "Normal"
$ws.Range("F8").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.we-online.com/katalog/datasheet/6941xx301002.pdf") | Out-Null
$ws.Range("F9").Style = "Normal"
$ws.Range("F9").Style = "Hyperlink"

# --- Number format for the whole price column (existing rows too) ---
$ws.Range("E3:E9").NumberFormat = "0.00"

# --- Column widths: B and C become one uniform width ---
$ws.Range("B:C").ColumnWidth = 22.140625

# --- Update selection to mirror the saved worksheet view ---
$ws.Range("C25").Select() | Out-Null

Write-Output "done"
